# Updated symbol list on Sun Dec 18 10:53:27 UTC 2022 with GitHub Actions
#
# Helper: write a numeric-looking value as TEXT (matching the sheet's
# existing inline-string cells) without leaving a stray NumberFormat/
# quote-prefix style behind on the cell.
function Set-TextValue {
    param($range, [string]$text)
    $range.Value = "'" + $text
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) refreshes --------------------------------------
Set-TextValue $ws.Range("D2")  "247.34"
Set-TextValue $ws.Range("D4")  "5.476"
Set-TextValue $ws.Range("D5")  "0.05619"
Set-TextValue $ws.Range("D6")  "6.469"
Set-TextValue $ws.Range("D7")  "0.8043"
Set-TextValue $ws.Range("D8")  "1.046"
Set-TextValue $ws.Range("D9")  "0.1426"
Set-TextValue $ws.Range("D10") "0.07263"
Set-TextValue $ws.Range("D12") "0.02971"
Set-TextValue $ws.Range("D13") "0.09261"
Set-TextValue $ws.Range("D14") "0.001676"
Set-TextValue $ws.Range("D15") "3.202"
Set-TextValue $ws.Range("D16") "0.04697"
Set-TextValue $ws.Range("D17") "0.0005982"
Set-TextValue $ws.Range("D18") "0.006285"
Set-TextValue $ws.Range("D19") "0.001056"
Set-TextValue $ws.Range("D20") "0.003819"
Set-TextValue $ws.Range("D25") "2.117"
Set-TextValue $ws.Range("D26") "0.3297"
Set-TextValue $ws.Range("D40") "0.04154"

# --- Rows 41-43 re-ranked: KickToken moves up to 41, BKEXToken to 42, --
# --- CEJI to 43 (each keeps its own refreshed price/volume label) -----
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.006842"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1043"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.002975"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Remaining price refreshes / label tweaks -------------------------
Set-TextValue $ws.Range("D44") "0.01029"
Set-TextValue $ws.Range("D45") "0.00005636"
Set-TextValue $ws.Range("D47") "0.6814"
Set-TextValue $ws.Range("D48") "0.02589"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
